# Swap columns B<->C and D<->E on the active sheet (column reorder refactor).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C").Cut()
$ws.Columns("B").Insert()

$ws.Columns("E").Cut()
$ws.Columns("D").Insert()

# Leave the same selection state as the final saved file (entire column D selected).
$ws.Range("D1:D1048576").Select()
